$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.668.78'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.68%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.846.93'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.20%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.013'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -2.92%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.95%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.011'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.74%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4309'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.85%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3750'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07351'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.80%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8808'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.59'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.852.61'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.730'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.456'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07121'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.63%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.61'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.38%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.013'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.97%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008981'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.010'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.50'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.676.46'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.259'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.17'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.76%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.085.72'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.038'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.65%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.54'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.59'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.139'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.390'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.75%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08932'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.89%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.229'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7767'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.559'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.62%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.903'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.90%  '
$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.142'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.49%  '
$ws.Range("B37").Value = 'Frax'
$ws.Range("C37").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.011'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.81%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05336'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.56%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01971'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.49%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.189'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.80%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.866'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5166'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.85%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1678'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.937'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '110.71'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.61%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.65'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4734'
$ws.Range("D47").Style = "Normal"
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06498'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.39%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.702'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.87%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.011'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.885'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.06%  '

Write-Output "Applied cryptos.xlsx update"